$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "TestSheet"
$ws.Range("A1").Value = "Hello"
$ws.Range("A2").Value = 42
$ws.Range("A3").Value = "123"
